$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kilavuzluk")

# Update header row: hizmet_turu -> tip, taban stays, ilave -> ek
$ws.Range("A1").Value = "tip"
$ws.Range("B1").Value = "taban"
$ws.Range("C1").Value = "ek"

# Add new data rows
$ws.Range("A2").Value = "bogaz_istanbul"
$ws.Range("B2").Value = 550
$ws.Range("C2").Value = 100

$ws.Range("A3").Value = "bogaz_canakkale"
$ws.Range("B3").Value = 550
$ws.Range("C3").Value = 100

$ws.Range("A4").Value = "halic"
$ws.Range("B4").Value = 605
$ws.Range("C4").Value = 136

# Make kilavuzluk the active / selected sheet
$ws.Select() | Out-Null
$ws.Range("A1:C4").Select() | Out-Null
